$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.392592191696167
$ws.Range("B1").Value = 2.623353719711304
$ws.Range("D1").Value = 1.483839511871338
$ws.Range("E1").Value = 0.9000283479690552
